$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2347844905132767
$ws.Range("C3").Value = 0.3440668887112966
$ws.Range("C4").Value = 0.4673456793719342
$ws.Range("C5").Value = 0.4252038582236517
$ws.Range("C6").Value = 0.394902939838619
